# The edit inserts one new data row into the table at sheet row 618,
# pushing the existing rows 618-703 down to 619-704 (dimension grows
# from A1:R703 to A1:R704), and fills the newly inserted row with its
# own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 618, shifting rows 618:703 down to 619:704.
$ws.Rows.Item(618).Insert()

# Populate the newly inserted row 618 with its data.
$ws.Cells.Item(618, 1).Value = 5
$ws.Cells.Item(618, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(618, 3).Value = "Maule"
$ws.Cells.Item(618, 4).Value = 45127
$ws.Cells.Item(618, 5).Value = 7
$ws.Cells.Item(618, 6).Value = 100112043
$ws.Cells.Item(618, 7).Value = "Pepino ensalada"
$ws.Cells.Item(618, 8).Value = "Sin especificar"
$ws.Cells.Item(618, 9).Value = "Primera"
$ws.Cells.Item(618, 10).Value = 300
$ws.Cells.Item(618, 11).Value = 9000
$ws.Cells.Item(618, 12).Value = 9000
$ws.Cells.Item(618, 13).Value = 9000
$ws.Cells.Item(618, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(618, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(618, 16).Value = 150
$ws.Cells.Item(618, 17).Value = 60
$ws.Cells.Item(618, 18).Value = "Hortaliza"
